# #5: cash & deposit done
# Finish out the 存款 (deposits) sheet: turn the stray "leftover header-style
# data row" in row 1 into a real header row, and append the same
# property_category/category/date/legislator_name/legislator_id/source_file/
# index metadata columns (G:M) that every other sheet in this workbook
# already carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- header row (row 1): B1:M1 -----------------------------------------
$headers = @{
    "B1" = "bank"
    "C1" = "deposit_type"
    "D1" = "currency"
    "E1" = "owner"
    "F1" = "total"
    "G1" = "property_category"
    "H1" = "category"
    "I1" = "date"
    "J1" = "legislator_name"
    "K1" = "legislator_id"
    "L1" = "source_file"
    "M1" = "index"
}
foreach ($addr in $headers.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $headers[$addr]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---- data rows 2-6 -------------------------------------------------------
# columns: A index | B bank | C deposit_type | D currency | E owner | F total
$rows = @(
    @{ Row=2; Index=62; Bank="臺灣銀行";         DepositType="活期存款"; Total=3343430 },
    @{ Row=3; Index=63; Bank="臺灣新光商業銀行"; DepositType="活期存款"; Total=896597 },
    @{ Row=4; Index=65; Bank="彰化商業銀行";     DepositType="活期存款"; Total=334231 },
    @{ Row=5; Index=66; Bank="屏東縣&山地區農會"; DepositType="活期存款"; Total=287527 },
    @{ Row=6; Index=67; Bank="中華郵政股份有限公司"; DepositType="活崩存款"; Total=4200 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Index          # A index
    $ws.Cells.Item($row, 2).Value = $r.Bank            # B bank
    $ws.Cells.Item($row, 3).Value = $r.DepositType     # C deposit_type
    $ws.Cells.Item($row, 4).Value = "新臺幣"           # D currency
    $ws.Cells.Item($row, 5).Value = "簡東明"           # E owner
    $ws.Cells.Item($row, 6).Value = $r.Total           # F total

    $ws.Cells.Item($row, 7).Value = "deposit"          # G property_category
    $ws.Cells.Item($row, 8).Value = "normal"           # H category
    $ws.Cells.Item($row, 9).Value = "2011-12-30"       # I date
    $ws.Cells.Item($row, 10).Value = "簡東明"          # J legislator_name
    $ws.Cells.Item($row, 11).Value = 1717              # K legislator_id
    $ws.Cells.Item($row, 12).Value = "tmp3d8a1"        # L source_file
    $ws.Cells.Item($row, 13).Value = $r.Index          # M index
}
